$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 173.73334
$ws.Range("I2").Value = 177.81818
$ws.Range("J2").Value = 162.5
$ws.Range("K2").Value = 177.81818
$ws.Range("L2").Value = 162.5
$ws.Range("M2").Value = -64.81818000000001
$ws.Range("N2").Value = -388.5
$ws.Range("H132").Value = 11192.952
$ws.Range("I132").Value = 11950.789
$ws.Range("K132").Value = 35852.367
$ws.Range("M132").Value = -33322.367
$ws.Range("H138").Value = 2824.75
$ws.Range("I138").Value = 2433.3333
$ws.Range("K138").Value = 7299.999899999999
$ws.Range("M138").Value = -2159.999899999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2152.1052
$ws.Range("I122").Value = 1661.5714
$ws.Range("J122").Value = 3525.6
$ws.Range("K122").Value = 4984.7142
$ws.Range("L122").Value = 10576.8
$ws.Range("M122").Value = -2534.7142
$ws.Range("N122").Value = -15476.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2003.75
$ws.Range("I68").Value = 1515
$ws.Range("K68").Value = 4545
$ws.Range("M68").Value = -3734
$ws.Range("H71").Value = 2003.75
$ws.Range("I71").Value = 1515
$ws.Range("K71").Value = 13635
$ws.Range("M71").Value = -9579
$ws.Range("H99").Value = 4965.6665
$ws.Range("I99").Value = 4965.6665
$ws.Range("K99").Value = 14896.9995
$ws.Range("M99").Value = -12650.9995
$ws.Range("H120").Value = 10449.5
$ws.Range("I120").Value = 10449.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 31348.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -26510.5
$ws.Range("H121").Value = 1490.6154
$ws.Range("I121").Value = 236.25
$ws.Range("J121").Value = 2048.111
$ws.Range("K121").Value = 708.75
$ws.Range("L121").Value = 6144.333
$ws.Range("M121").Value = 601.25
$ws.Range("N121").Value = -8764.332999999999
$ws.Range("H122").Value = 990
$ws.Range("I122").Value = 990
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8910
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6460
$ws.Range("H123").Value = 799
$ws.Range("I123").Value = 799
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 2397
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = 53
$ws.Range("H124").Value = 900
$ws.Range("I124").Value = 900
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 2700
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 2210
$ws.Range("H125").Value = 6965
$ws.Range("I125").Value = 3930
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 11790
$ws.Range("L125").Value = 30000
$ws.Range("M125").Value = -6870
$ws.Range("N125").Value = -39840
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 275994.4
$ws.Range("I128").Value = 275994.4
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 827983.2000000001
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -823003.2000000001
$ws.Range("H129").Value = 1916
$ws.Range("I129").Value = 450
$ws.Range("J129").Value = 2282.5
$ws.Range("K129").Value = 1350
$ws.Range("L129").Value = 6847.5
$ws.Range("M129").Value = 3650
$ws.Range("N129").Value = -16847.5
$ws.Range("H130").Value = 7872.6665
$ws.Range("I130").Value = 585
$ws.Range("J130").Value = 11516.5
$ws.Range("K130").Value = 1755
$ws.Range("L130").Value = 34549.5
$ws.Range("M130").Value = 3265
$ws.Range("N130").Value = -44589.5
$ws.Range("H131").Value = 2505.8333
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2505.8333
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 7517.499899999999
$ws.Range("N131").Value = -17597.4999
$ws.Range("H132").Value = 1991.4166
$ws.Range("I132").Value = 1739.7
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 15657.3
$ws.Range("L132").Value = 29250
$ws.Range("M132").Value = -13127.3
$ws.Range("N132").Value = -34310
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 2015
$ws.Range("I134").Value = 2015
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6045
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -975
$ws.Range("H136").Value = 7250.4287
$ws.Range("I136").Value = 4499
$ws.Range("J136").Value = 8351
$ws.Range("K136").Value = 13497
$ws.Range("L136").Value = 25053
$ws.Range("M136").Value = -8397
$ws.Range("N136").Value = -35253
$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 3000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3900
$ws.Range("H138").Value = 3694.818
$ws.Range("I138").Value = 2613.2856
$ws.Range("J138").Value = 5587.5
$ws.Range("K138").Value = 7839.8568
$ws.Range("L138").Value = 16762.5
$ws.Range("M138").Value = -2699.8568
$ws.Range("N138").Value = -27042.5
$ws.Range("H139").Value = 3371.5
$ws.Range("I139").Value = 2628.6667
$ws.Range("J139").Value = 5600
$ws.Range("K139").Value = 7886.000100000001
$ws.Range("L139").Value = 16800
$ws.Range("M139").Value = -2746.000100000001
$ws.Range("N139").Value = -27080
$ws.Range("H140").Value = 2321.1428
$ws.Range("I140").Value = 2041.3334
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 6124.0002
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -944.0002000000004
$ws.Range("N140").Value = -22360
$ws.Range("H141").Value = 4276.3335
$ws.Range("I141").Value = 1414.5
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 4243.5
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = 936.5
$ws.Range("N141").Value = -40360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 105.15385
$ws.Range("I2").Value = 63.444443
$ws.Range("K2").Value = 63.444443
$ws.Range("M2").Value = 49.555557
$ws.Range("H97").Value = 949.75
$ws.Range("H102").Value = 1397.6923
$ws.Range("I102").Value = 1397.6923
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1397.6923
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 224.3077000000001
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1182.1538
$ws.Range("I122").Value = 1199.8334
$ws.Range("K122").Value = 3599.5002
$ws.Range("M122").Value = -1149.5002
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 26000
$ws.Range("J121").Value = 26000
$ws.Range("L121").Value = 26000
$ws.Range("N121").Value = -29494
$ws.Range("H122").Value = 2492.5715
$ws.Range("I122").Value = 2289.8
$ws.Range("K122").Value = 6869.400000000001
$ws.Range("M122").Value = -4419.400000000001
$ws.Range("H132").Value = 3499.75
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 3022.3
$ws.Range("I136").Value = 2358.111
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 7074.333
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -4524.333
$ws.Range("N136").Value = -32100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17144.385
$ws.Range("J2").Value = 11363.857
$ws.Range("L2").Value = 11363.857
$ws.Range("N2").Value = -11587.857
$ws.Range("H40").Value = 39999.5
$ws.Range("J40").Value = 19999
$ws.Range("L40").Value = 19999
$ws.Range("N40").Value = -20297
$ws.Range("H81").Value = 11333
$ws.Range("J81").Value = 19999.5
$ws.Range("L81").Value = 39999
$ws.Range("N81").Value = -42121
$ws.Range("H84").Value = 11333
$ws.Range("J84").Value = 19999.5
$ws.Range("L84").Value = 199995
$ws.Range("N84").Value = -210603
$ws.Range("H122").Value = 4748.8335
$ws.Range("I122").Value = 4726
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14178
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -11728
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2633.8333
$ws.Range("I132").Value = 2101.3333
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 6303.999899999999
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -3773.999899999999
$ws.Range("N132").Value = -14558.9999
$ws.Range("H136").Value = 3056.5454
$ws.Range("I136").Value = 1997.0454
$ws.Range("J136").Value = 5175.5454
$ws.Range("K136").Value = 5991.1362
$ws.Range("L136").Value = 15526.6362
$ws.Range("M136").Value = -3441.1362
$ws.Range("N136").Value = -20626.6362
